# "gpio and lcd descriptions added"
#
# Adds a "Примечание" (Note) column to the power-consumption table on the
# active sheet:
#   - F1 gets a new header cell "Примечание", styled like the other
#     header cells (copy formatting from B1).
#   - D10 (TDA2003V current draw, previously a placeholder "??") is set
#     to the real value of 1 (A).
#   - F10 gets the explanatory note about the TDA2003V power budget.
#   - Column F is widened so the long note text fits.
#   - The active selection is moved to F18 to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Примечание" header in F1, formatted like the other headers ---
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Value = "Примечание"

# --- TDA2003V row: replace the "??" placeholder current with the real value ---
$ws.Range("D10").Value = 1

# --- Explanation note for the TDA2003V power/current calculation ---
$ws.Range("F10").Value = "U вых max = 6 В, U вых сред.кв = 5 В, R = 4 Ом => P = 6 Вт . КПД ~= 50%, потребление около 12 Вт => ток 1 А по 14.4 В"

# --- Widen column F so the note is readable ---
$ws.Columns.Item(6).ColumnWidth = 126.42

# --- Match the saved cursor/selection position ---
[void]$ws.Range("F18").Select()

Write-Host "Updated power consumption sheet: GPIO/LCD notes added"
